# Applies the scheduled-runner price/profit refresh described in the commit
# diff. Each sheet's rows keep their Leve metadata (A-G) untouched; only the
# market-price/profit columns H-N are refreshed. Where the diff shows a cell
# disappearing entirely (no replacement <c> element), we ClearContents() so
# the cell is dropped from the saved XML instead of being written as 0/blank.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------- ALC sheet
$ws = $wb.Worksheets.Item("ALC")

$ws.Range("H43").Value = 3247.5
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 3247.5
$ws.Range("K43").Value = 0
$ws.Range("L43").ClearContents()
$ws.Range("M43").Value = 3247.5
$ws.Range("N43").Value = -3385.5

$ws.Range("H69").Value = 0
$ws.Range("J69").Value = 0
$ws.Range("L69").ClearContents()
$ws.Range("N69").Value = 0

$ws.Range("H72").Value = 0
$ws.Range("J72").Value = 0
$ws.Range("L72").ClearContents()
$ws.Range("N72").Value = 0

$ws.Range("H80").Value = 457.42856
$ws.Range("I80").Value = 125.5
$ws.Range("J80").Value = 900
$ws.Range("K80").Value = 376.5
$ws.Range("L80").Value = 2700
$ws.Range("M80").Value = 621.5
$ws.Range("N80").Value = -4696

$ws.Range("H83").Value = 457.42856
$ws.Range("I83").Value = 125.5
$ws.Range("J83").Value = 900
$ws.Range("K83").Value = 1129.5
$ws.Range("L83").Value = 8100
$ws.Range("M83").Value = 3862.5
$ws.Range("N83").Value = -18084

$ws.Range("H138").Value = 11714
$ws.Range("J138").Value = 12999.667
$ws.Range("L138").Value = 38999.001
$ws.Range("N138").Value = -49279.001

# ---------------------------------------------------------------- BSM sheet
$ws = $wb.Worksheets.Item("BSM")

$ws.Range("H18").Value = 6999
$ws.Range("J18").Value = 6999
$ws.Range("L18").Value = 6999
$ws.Range("N18").Value = -8057

$ws.Range("H23").Value = 800
$ws.Range("I23").Value = 800
$ws.Range("K23").Value = 800
$ws.Range("M23").Value = -517

$ws.Range("H46").Value = 0
$ws.Range("I46").Value = 0
$ws.Range("K46").Value = 0
$ws.Range("M46").ClearContents()

$ws.Range("H54").Value = 5333.3335
$ws.Range("I54").Value = 5333.3335
$ws.Range("K54").Value = 5333.3335
$ws.Range("M54").Value = -4849.3335

$ws.Range("H56").Value = 5000
$ws.Range("I56").Value = 5000
$ws.Range("K56").Value = 5000
$ws.Range("M56").Value = -4261

$ws.Range("H86").Value = 2101.2
$ws.Range("I86").Value = 2166.3333
$ws.Range("K86").Value = 2166.3333
$ws.Range("M86").Value = -1043.3333

$ws.Range("H89").Value = 2101.2
$ws.Range("I89").Value = 2166.3333
$ws.Range("K89").Value = 10831.6665
$ws.Range("M89").Value = -5215.666499999999

$ws.Range("H94").Value = 2499.3333
$ws.Range("I94").Value = 2874.75
$ws.Range("J94").Value = 1748.5
$ws.Range("K94").Value = 2874.75
$ws.Range("L94").Value = 1748.5
$ws.Range("M94").Value = -2423.75
$ws.Range("N94").Value = -2650.5

$ws.Range("H134").Value = 6900
$ws.Range("I134").Value = 4450
$ws.Range("K134").Value = 13350
$ws.Range("M134").Value = -10815

$ws.Range("H140").Value = 77499.25
$ws.Range("J140").Value = 77499.25
$ws.Range("L140").Value = 77499.25
$ws.Range("N140").Value = -87859.25

# ---------------------------------------------------------------- CRP sheet
$ws = $wb.Worksheets.Item("CRP")

$ws.Range("H19").Value = 288.16666
$ws.Range("I19").Value = 155.8
$ws.Range("J19").Value = 950
$ws.Range("K19").Value = 155.8
$ws.Range("L19").Value = 950
$ws.Range("M19").Value = 14.19999999999999
$ws.Range("N19").Value = -1290

$ws.Range("H24").Value = 288.16666
$ws.Range("I24").Value = 155.8
$ws.Range("J24").Value = 950
$ws.Range("K24").Value = 155.8
$ws.Range("L24").Value = 950
$ws.Range("M24").Value = 14.19999999999999
$ws.Range("N24").Value = -1290

$ws.Range("H39").Value = 0
$ws.Range("I39").Value = 0
$ws.Range("K39").Value = 0
$ws.Range("M39").ClearContents()

$ws.Range("H49").Value = 0
$ws.Range("I49").Value = 0
$ws.Range("K49").Value = 0
$ws.Range("M49").ClearContents()

$ws.Range("H62").Value = 0
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("L62").ClearContents()
$ws.Range("M62").ClearContents()
$ws.Range("N62").Value = 0

$ws.Range("H65").Value = 0
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("L65").ClearContents()
$ws.Range("M65").ClearContents()
$ws.Range("N65").Value = 0

$ws.Range("H99").Value = 0
$ws.Range("I99").Value = 0
$ws.Range("K99").Value = 0
$ws.Range("M99").ClearContents()

$ws.Range("H126").Value = 0
$ws.Range("I126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("M126").ClearContents()

$ws.Range("H141").Value = 79249.5
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 79249.5
$ws.Range("K141").Value = 0
$ws.Range("L141").ClearContents()
$ws.Range("M141").Value = 79249.5
$ws.Range("N141").Value = -89609.5

# ---------------------------------------------------------------- CUL sheet
$ws = $wb.Worksheets.Item("CUL")

$ws.Range("H139").Value = 50000
$ws.Range("I139").Value = 50000
$ws.Range("K139").Value = 150000
$ws.Range("M139").Value = -144860

# ---------------------------------------------------------------- GSM sheet
$ws = $wb.Worksheets.Item("GSM")

$ws.Range("H118").Value = 16097
$ws.Range("J118").Value = 16097
$ws.Range("L118").Value = 16097
$ws.Range("N118").Value = -19411

# ---------------------------------------------------------------- LTW sheet
$ws = $wb.Worksheets.Item("LTW")

$ws.Range("H68").Value = 5549.25
$ws.Range("I68").Value = 5065.6665
$ws.Range("J68").Value = 7000
$ws.Range("K68").Value = 5065.6665
$ws.Range("L68").Value = 7000
$ws.Range("M68").Value = -4316.6665
$ws.Range("N68").Value = -8498

$ws.Range("H71").Value = 5549.25
$ws.Range("I71").Value = 5065.6665
$ws.Range("J71").Value = 7000
$ws.Range("K71").Value = 25328.3325
$ws.Range("L71").Value = 35000
$ws.Range("M71").Value = -21584.3325
$ws.Range("N71").Value = -42488

$ws.Range("H82").Value = 2037.625
$ws.Range("I82").Value = 908.25
$ws.Range("J82").Value = 3167
$ws.Range("K82").Value = 908.25
$ws.Range("L82").Value = 3167
$ws.Range("M82").Value = -547.25
$ws.Range("N82").Value = -3889

$ws.Range("H85").Value = 2037.625
$ws.Range("I85").Value = 908.25
$ws.Range("J85").Value = 3167
$ws.Range("K85").Value = 908.25
$ws.Range("L85").Value = 3167
$ws.Range("M85").Value = 339.75
$ws.Range("N85").Value = -5663

$ws.Range("H104").Value = 14273.8
$ws.Range("J104").Value = 14273.8
$ws.Range("L104").Value = 14273.8
$ws.Range("N104").Value = -21261.8
